$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- B11: high-pass ratio 0.05 -> 0.1 ---
$ws.Range("B11").Value = 0.1

# --- E16: "文件时间是否为UTC" False -> True ---
# Use Copy + PasteSpecial(values) from the Sheet2 "True" list source so the
# cell keeps its shared-string/text representation instead of being coerced
# into a native Excel boolean (xlPasteValues = -4163).
$ws2.Range("B1").Copy()
$ws.Range("E16").PasteSpecial(-4163)

# --- B23: "是否修正：平场" False -> True (same technique) ---
$ws2.Range("B1").Copy()
$ws.Range("B23").PasteSpecial(-4163)

# --- New row 25/26 cells: flat-field channel list settings ---
# Write values first (in the same order the target workbook's shared
# strings were introduced), then copy the D23/E23 formatting onto them
# (xlPasteFormats = -4122) so the new cells match styles s=7 / s=10.
$ws.Range("E25").Value = "HSOLRGB"
$ws.Range("D25").Value = "平场通道列表"
$ws.Range("D23").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E25").PasteSpecial(-4122)

$ws.Range("E26").Value = "_Bin2_"
$ws.Range("D26").Value = "通道名格式（紧邻前缀）"
$ws.Range("D23").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E23").Copy()
$ws.Range("E26").PasteSpecial(-4122)

# --- B25/B26/B27: point bias/dark/flat paths at the new bias-dark-flat tree ---
$ws.Range("B25").Value = "/home/hao/astro/bias-dark-flat/294mm-pro-bin2/bias-master.fits"
$ws.Range("B26").Value = "/home/hao/astro/bias-dark-flat/294mm-pro-bin2/dark-master.fits"
$ws.Range("B27").Value = "/home/hao/astro/bias-dark-flat"

# --- Move the active selection like the saved workbook ---
$ws.Range("E29").Select() | Out-Null
$excel.CutCopyMode = $false
